# Apply crypto price/volume updates (and two row swaps) per the commit diff.
# Values that look like plain numbers are forced to remain text (matching the
# original inline-string cell content, e.g. "1.00" must not become the number 1),
# by temporarily marking the cell as Text, assigning the value, then clearing the
# format again so the cell keeps its original (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.111.58'
$ws.Range('E2').Value = '  -6.16%  '
$ws.Range('D3').Value = '2.443.45'
$ws.Range('E3').Value = '  -8.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '530.70'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -3.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.04'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.565'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -4.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0984'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -6.75%  '
$ws.Range('E10').Value = '  -2.59%  '
$ws.Range('E11').Value = '  +4.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.348'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -5.52%  '
$ws.Range('D13').Value = '2.873.89'
$ws.Range('E13').Value = '  -8.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '24.00'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -7.64%  '
$ws.Range('D15').Value = '59.100.25'
$ws.Range('E15').Value = '  -5.96%  '
$ws.Range('E16').Value = '  -7.07%  '
$ws.Range('D17').Value = '2.494.29'
$ws.Range('E17').Value = '  -7.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.07'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -7.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.35'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -5.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '323.55'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -5.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.969'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.00%  '
$ws.Range('E22').Value = '  -9.07%  '
$ws.Range('E23').Value = '  -7.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.00'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -5.50%  '
$ws.Range('E25').Value = '  -4.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.974'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.68'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -5.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.29'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -3.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.85'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.43%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.81'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -5.95%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0770'
$ws.Range('E31').Value = '  -9.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.998'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '157.66'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.50'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -6.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.23'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.60%  '
$ws.Range('E36').Value = '  -6.22%  '
$ws.Range('E37').Value = '  -2.81%  '
$ws.Range('E38').Value = '  -6.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '312.29'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -7.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.853'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -8.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.71'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.70'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -6.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.995'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('E44').Value = '  -3.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0933'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -4.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0520'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -7.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.98'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -8.66%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0228'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.45'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -9.24%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '1.980.93'
$ws.Range('E51').Value = '  -5.20%  '
